$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder fabric names with the real ones
$ws.Range("A2").Value = "ALGODÓN"
$ws.Range("A3").Value = "NYLON"
$ws.Range("A4").Value = "SEDA"
$ws.Range("A5").Value = "POLIESTER"

# The old duplicate row is now left blank
$ws.Range("A6").ClearContents()

# A2 gets its own look (Calibri) and is centered
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").HorizontalAlignment = -4108

# A3:A6 are centered, Arial 11pt
$ws.Range("A3:A6").Font.Name = "Arial"
$ws.Range("A3:A6").Font.Size = 11
$ws.Range("A3:A6").HorizontalAlignment = -4108

# Column A is widened to fit the longer fabric names
$ws.Columns.Item(1).ColumnWidth = 71.6

# Selection ends up on A6 (the now-empty row)
$ws.Range("A6").Select() | Out-Null
